$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Introduction to R"
$ws.Range("H2").Value = 5

# Size the new column the way Excel would after a double-click
# "AutoFit" on the column border (matches the bestFit width Excel
# stored for this header text).
$ws.Columns.Item(8).ColumnWidth = 15.14
